# Insert a new daily-price record for Mango (Vega Central Mapocho de Santiago)
# at row 198, shifting the existing rows 198:330 down to 199:331.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 198 (pushes old row198..row330 to 199..331)
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new observation.
$ws.Range("A198").Value = 9
$ws.Range("B198").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C198").Value = "Metropolitana"
$ws.Range("D198").Value = 44574
$ws.Range("E198").Value = 13
$ws.Range("F198").Value = "Fruta"
$ws.Range("G198").Value = 100108
$ws.Range("H198").Value = "Tropicales y subtropicales"
$ws.Range("I198").Value = 100108002
$ws.Range("J198").Value = "Mango"
$ws.Range("K198").Value = "Sin especificar"
$ws.Range("L198").Value = "Primera"
$ws.Range("M198").Value = 580
$ws.Range("N198").Value = 6000
$ws.Range("O198").Value = 6000
$ws.Range("P198").Value = 6000
$ws.Range("Q198").Value = "`$/bandeja 4 kilos"
$ws.Range("R198").Value = "Perú"
$ws.Range("S198").Value = 1500
$ws.Range("T198").Value = 4
